$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text storage for numeric-looking price strings so Excel
# does not silently convert them to numbers and strip formatting/precision.
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

# Apply updated cell values
$ws.Range('D2').Value = '64.738.62'
$ws.Range('E2').Value = '  -0.84%  '
$ws.Range('D3').Value = '3.136.35'
$ws.Range('E3').Value = '  -1.31%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '579.99'
$ws.Range('E5').Value = '  +0.92%  '
$ws.Range('D6').Value = '146.98'
$ws.Range('E6').Value = '  -2.73%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('D8').Value = '3.134.94'
$ws.Range('E8').Value = '  -1.22%  '
$ws.Range('D9').Value = '0.524'
$ws.Range('E9').Value = '  -0.95%  '
$ws.Range('E10').Value = '  -3.96%  '
$ws.Range('E11').Value = '  -2.34%  '
$ws.Range('D12').Value = '0.496'
$ws.Range('E12').Value = '  -2.96%  '
$ws.Range('E13').Value = '  -2.54%  '
$ws.Range('D14').Value = '36.98'
$ws.Range('E14').Value = '  -3.82%  '
$ws.Range('D15').Value = '3.651.73'
$ws.Range('E15').Value = '  -0.82%  '
$ws.Range('D16').Value = '64.793.42'
$ws.Range('E16').Value = '  -0.75%  '
$ws.Range('D17').Value = '3.139.31'
$ws.Range('E17').Value = '  -0.85%  '
$ws.Range('D18').Value = '7.12'
$ws.Range('E18').Value = '  -2.08%  '
$ws.Range('E19').Value = '  -0.76%  '
$ws.Range('D20').Value = '498.17'
$ws.Range('E20').Value = '  -3.36%  '
$ws.Range('D21').Value = '15.25'
$ws.Range('E21').Value = '  +1.36%  '
$ws.Range('D22').Value = '0.710'
$ws.Range('E22').Value = '  -4.49%  '
$ws.Range('D23').Value = '14.98'
$ws.Range('E23').Value = '  -7.66%  '
$ws.Range('D24').Value = '7.72'
$ws.Range('E24').Value = '  -2.38%  '
$ws.Range('D25').Value = '84.14'
$ws.Range('E25').Value = '  -1.20%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('E27').Value = '  -1.70%  '
$ws.Range('E28').Value = '  -0.41%  '
$ws.Range('E29').Value = '  -1.87%  '
$ws.Range('E30').Value = '  +0.68%  '
$ws.Range('D31').Value = '27.40'
$ws.Range('E31').Value = '  -2.81%  '
$ws.Range('E32').Value = '  -1.49%  '
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('D34').Value = '6.39'
$ws.Range('E34').Value = '  +0.95%  '
$ws.Range('D35').Value = '6.43'
$ws.Range('E35').Value = '  -3.85%  '
$ws.Range('D36').Value = '54.76'
$ws.Range('E36').Value = '  -2.17%  '
$ws.Range('D37').Value = '0.0889'
$ws.Range('E37').Value = '  +1.14%  '
$ws.Range('D38').Value = '467.57'
$ws.Range('E38').Value = '  -2.38%  '
$ws.Range('D39').Value = '0.0416'
$ws.Range('E39').Value = '  -1.58%  '
$ws.Range('E40').Value = '  -6.44%  '
$ws.Range('D41').Value = '8.70'
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('D42').Value = '2.968.88'
$ws.Range('E42').Value = '  -5.32%  '
$ws.Range('E43').Value = '  -4.54%  '
$ws.Range('E44').Value = '  -4.63%  '
$ws.Range('D45').Value = '0.281'
$ws.Range('E45').Value = '  -4.13%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').Value = '28.14'
$ws.Range('E46').Value = '  -4.44%  '
$ws.Range('B47').Value = 'PEPE'
$ws.Range('C47').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D47').Value = '0.0₃0596'
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('E49').Value = '  -2.34%  '
$ws.Range('D50').Value = '2.22'
$ws.Range('E50').Value = '  -5.05%  '
$ws.Range('D51').Value = '119.12'
$ws.Range('E51').Value = '  -4.51%  '
